$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegistrationPageTests")

# New testcase row: duplicate the existing data row (2) into row 3, keeping
# the same FirstName/LastName/Mobile/Email/Password sample data, then
# rename the testcase identifier in column A for the new
# "two users, same email" duplicate-email scenario.
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4104) | Out-Null

$ws.Range("A3").Value = "add_TwoUsers_With_Same_Email"

# Re-create the mailto hyperlink on the new row's Email cell, and restore
# the "Hyperlink" cell style (PasteSpecial only copied the plain value).
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:a@gmail.com") | Out-Null
$ws.Range("E3").Style = "Hyperlink"

# The new, longer testcase name no longer fits column A - widen it.
$ws.Columns.Item(1).AutoFit() | Out-Null
